$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.644.28"
$ws.Cells.Item(2, 5).Value = "  -0.97%  "
$ws.Cells.Item(3, 4).Value = "2.434.99"
$ws.Cells.Item(3, 5).Value = "  -1.49%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "506.75"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.57%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "128.78"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.90%  "
$ws.Cells.Item(7, 5).Value = "  -0.07%  "
$ws.Cells.Item(8, 5).Value = "  -1.39%  "
$ws.Cells.Item(9, 4).Value = "2.447.67"
$ws.Cells.Item(9, 5).Value = "  -1.10%  "
$ws.Cells.Item(10, 5).Value = "  -0.34%  "
$ws.Cells.Item(11, 5).Value = "  -4.13%  "
$ws.Cells.Item(12, 5).Value = "  -3.50%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.331"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -3.43%  "
$ws.Cells.Item(14, 4).Value = "2.867.65"
$ws.Cells.Item(14, 5).Value = "  -1.44%  "
$ws.Cells.Item(15, 4).Value = "57.569.42"
$ws.Cells.Item(15, 5).Value = "  -0.96%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "21.79"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.62%  "
$ws.Cells.Item(17, 5).Value = "  -2.93%  "
$ws.Cells.Item(18, 4).Value = "2.444.96"
$ws.Cells.Item(18, 5).Value = "  -1.18%  "
$ws.Cells.Item(19, 5).Value = "  -3.98%  "
$ws.Cells.Item(20, 5).Value = "  -1.83%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "314.88"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.03%  "
$ws.Cells.Item(22, 5).Value = "  -0.13%  "
$ws.Cells.Item(23, 5).Value = "  -1.42%  "
$ws.Cells.Item(24, 5).Value = "  -1.64%  "
$ws.Cells.Item(25, 5).Value = "  -0.43%  "
$ws.Cells.Item(27, 5).Value = "  -1.07%  "
$ws.Cells.Item(28, 5).Value = "  -2.58%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "169.77"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.37%  "
$ws.Cells.Item(30, 5).Value = "  -3.85%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.20"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.06%  "
$ws.Cells.Item(33, 5).Value = "  +0.23%  "
$ws.Cells.Item(34, 5).Value = "  -0.04%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.996"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.23%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "17.72"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -2.40%  "
$ws.Cells.Item(37, 5).Value = "  -5.53%  "
$ws.Cells.Item(38, 5).Value = "  -2.25%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "36.29"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.58%  "
$ws.Cells.Item(40, 5).Value = "  -2.63%  "
$ws.Cells.Item(41, 5).Value = "  -4.83%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "270.86"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.43%  "
$ws.Cells.Item(43, 5).Value = "  -2.96%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "4.87"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.32%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.580"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -2.85%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0909"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.05%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "119.65"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -5.42%  "
$ws.Cells.Item(48, 5).Value = "  -1.74%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "17.09"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -4.38%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0210"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.15%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "16.60"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.33%  "
